$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.588.45'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.919.14'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.11'
$ws.Range("E5").Value = '  -0.85%  '
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4860'
$ws.Range("E7").Value = '  +3.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06733'
$ws.Range("E9").Value = '  -1.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '111.74'
$ws.Range("E10").Value = '  +6.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.37'
$ws.Range("E11").Value = '  +5.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.918.70'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07589'
$ws.Range("E13").Value = '  -1.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.383'
$ws.Range("E14").Value = '  +1.92%  '
$ws.Range("E15").Value = '  +0.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '295.66'
$ws.Range("E16").Value = '  +2.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.583.28'
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.08'
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007566'
$ws.Range("E20").Value = '  -0.41%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.524'
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.171.63'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.428'
$ws.Range("E24").Value = '  +1.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.493'
$ws.Range("E25").Value = '  +1.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.97'
$ws.Range("E26").Value = '  -2.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.43'
$ws.Range("E27").Value = '  -3.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.108'
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("E29").Value = '  +0.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.434'
$ws.Range("E30").Value = '  +2.84%  '
$ws.Range("E31").Value = '  -0.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.103'
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05017'
$ws.Range("E33").Value = '  -0.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7413'
$ws.Range("E34").Value = '  +0.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.140'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9996'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02027'
$ws.Range("E37").Value = '  -2.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.703'
$ws.Range("E38").Value = '  -1.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.692'
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.022'
$ws.Range("E40").Value = '  -1.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '109.84'
$ws.Range("E41").Value = '  -1.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4455'
$ws.Range("E42").Value = '  +1.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8666'
$ws.Range("E43").Value = '  -1.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.885'
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.95'
$ws.Range("E45").Value = '  +4.18%  '
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.266'
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.309'
$ws.Range("E48").Value = '  +0.48%  '
$ws.Range("B49").Value = 'BitcoinSV'
$ws.Range("C49").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '48.39'
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("B50").Value = 'WOONetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.2551'
$ws.Range("E50").Value = '  +2.77%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1230'
$ws.Range("E51").Value = '  -0.15%  '
